$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# --- Add new column C mapping formula for rows 2-45 ---
# Row 2 gets its own (non-shared) formula; rows 3:45 are filled as a shared formula block.
$ws.Range("C2").Formula = '=_xlfn.CONCAT(A2," as ",B2, ",")'
$ws.Range("C3:C45").Formula = '=_xlfn.CONCAT(A3," as ",B3, ",")'

# Filling the formula into C3:C45 should not carry over the pre-existing cell
# formatting (yellow fill on C33:C45, and the blank-but-styled C29). Reset those
# cells back to the workbook's default "Normal" style, matching a plain fill-down
# from the unstyled C2 cell.
$ws.Range("C3:C45").Style = "Normal"

# --- Clear the old "column" / "group" header labels in row 32 ---
$ws.Range("B32:C32").ClearContents()

# --- Update sheet view: selection moved to B11, scrolled back to show column A ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()
$wb.Save()
